# LMS-2340 Updating basynthec stuff based on results of the leiden meeting.
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("openbis-metadata")
$ws2 = $wb.Worksheets.Item("openbis-data")

# --- match the saved workbook window geometry ---
$win = $wb.Windows.Item(1)
$win.Top = -20
$win.Left = -20
$win.Width = 24720
$win.Height = 16740

# --- openbis-metadata sheet: remove the "Strain" property row (was row 3) ---
$ws1.Rows.Item(3).Delete() | Out-Null

# --- openbis-data sheet: rename header + add additional strain rows ---
$ws2.Range("A1").Value = "Strain"
$ws2.Range("A2").Value = "MGP1"

$dataRow = $ws2.Range("C2:U2").Value()

$ws2.Range("A3").Value = "MGP100"
$ws2.Range("B3").Value = "OD600"
$ws2.Range("C3:U3").Value = $dataRow

$ws2.Range("A4").Value = "MGP20"
$ws2.Range("B4").Value = "OD600"
$ws2.Range("C4:U4").Value = $dataRow

$ws2.Range("A5").Value = "MGP999"
$ws2.Range("B5").Value = "OD600"
$ws2.Range("C5:U5").Value = $dataRow

# --- restore the selections shown by each sheet when the file was saved ---
$ws1.Rows.Item(3).EntireRow.Select() | Out-Null
$ws2.Range("A13").Select() | Out-Null
$ws1.Activate()
